$d = $word.ActiveDocument

# New bullet points to append after "Complete game and generate results",
# describing the Cover Story System design notes. "level" is the Word
# ListFormat.ListLevelNumber (1-based; ilvl 0 -> level 1, ilvl 1 -> level 2).
$items = @(
    @{ text = "Cover Story System"; level = 1 },
    @{ text = "Interacting with the date creates a “cover story”"; level = 2 },
    @{ text = "With each dialogue choice you pick, you have to make sure that the rest of your choices line up with what you stated previously"; level = 2 },
    @{ text = "If a dialogue option doesn’t add up, the date will grow suspicious, lowering the date quality score"; level = 2 },
    @{ text = "If the date becomes too suspicious, the date scenario will fail and the player will lose"; level = 2 }
)

foreach ($item in $items) {
    $lastPara = $d.Paragraphs.Last
    $r = $lastPara.Range
    $r.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    $newPara.Range.Text = $item.text
    $newPara.Range.ListFormat.ListLevelNumber = $item.level
}

# The document carries a collapsed "_GoBack" bookmark marking the last edit
# position. Originally it sat right after "Complete game and generate
# results"; since that is where the new text was typed, move it to the end
# of the newly added content (after "...player will lose").
#
# A collapsed Range positioned exactly at a paragraph's end boundary is
# mishandled when handed to Bookmarks.Add, so: append a throwaway trailing
# paragraph with placeholder text, bookmark the (unambiguous) start of that
# paragraph under the name "_GoBack" -- which both relocates the existing
# bookmark and removes it from its old spot -- delete the placeholder text,
# then delete the paragraph mark that separated it from the real last
# paragraph so everything merges back into a single paragraph with the
# bookmark collapsed right after "lose".
$trailPara = $d.Paragraphs.Last
$trailPara.Range.InsertParagraphAfter()
$trailNew = $d.Paragraphs.Last
$trailNew.Range.Text = "ZZZPLACEHOLDERZZZ"

$trailNew2 = $d.Paragraphs.Last
$targetPos = $trailNew2.Range.Start
$bmRange = $d.Range($targetPos, $targetPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$trailNew3 = $d.Paragraphs.Last
$textDelRange = $d.Range($targetPos, $trailNew3.Range.End - 1)
$textDelRange.Delete()

$joinRange = $d.Range($targetPos - 1, $targetPos)
$joinRange.Delete()
